# EOD report: append the 02 Jan 2019 closing entries (6 rows) below the
# existing 27 Dec 2018 data, mirroring the existing date/amount layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$date = "02 Jan 2019"
$amount = 499.8

$firstRow = 50
$lastRow = 55

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 2).Value = $amount
}
